$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sliding-window results were regenerated after refactoring the code that
# saves results into a specified output folder; the underlying predictions
# (and thus IPC PO / DELTA / DELTA^2, and the TOTAL/SSE/MSE summary) changed.
$ws.Range("C2").Value = 29.50128534542506
$ws.Range("D2").Value = -0.3387146545749395
$ws.Range("E2").Value = 0.1147276172238206
$ws.Range("C3").Value = 29.72259393625888
$ws.Range("D3").Value = -0.08740606374112048
$ws.Range("E3").Value = 0.007639819978716816
$ws.Range("C4").Value = 30.18642443566212
$ws.Range("D4").Value = 0.2664244356621168
$ws.Range("E4").Value = 0.07098197991787743
$ws.Range("C5").Value = 29.7076232420457
$ws.Range("D5").Value = -0.272376757954305
$ws.Range("E5").Value = 0.07418909827369807
$ws.Range("C6").Value = 29.72710782357337
$ws.Range("D6").Value = -0.312892176426633
$ws.Range("E6").Value = 0.09790151406899525
$ws.Range("C7").Value = 29.46802337815439
$ws.Range("D7").Value = -0.741976621845609
$ws.Range("E7").Value = 0.5505293073654218
$ws.Range("C8").Value = 29.37873299744113
$ws.Range("D8").Value = -0.8412670025588724
$ws.Range("E8").Value = 0.7077301695943898
$ws.Range("C9").Value = 29.97453232663057
$ws.Range("D9").Value = -0.4054676733694258
$ws.Range("E9").Value = 0.1644040341476154
$ws.Range("C10").Value = 30.16820205594819
$ws.Range("D10").Value = -0.2717979440518121
$ws.Range("E10").Value = 0.073874122390792
$ws.Range("C11").Value = 30.30205841598819
$ws.Range("D11").Value = -0.1779415840118119
$ws.Range("E11").Value = 0.03166320732063272
$ws.Range("C12").Value = 30.92585986026115
$ws.Range("D12").Value = 0.2358598602611508
$ws.Range("E12").Value = 0.05562987368240958
$ws.Range("C13").Value = 30.68201219889807
$ws.Range("D13").Value = -0.06798780110193192
$ws.Range("E13").Value = 0.004622341098675855
$ws.Range("C14").Value = 31.4800043214025
$ws.Range("D14").Value = 0.540004321402499
$ws.Range("E14").Value = 0.2916046671333734
$ws.Range("C15").Value = 31.10581343353556
$ws.Range("D15").Value = 0.1558134335355561
$ws.Range("E15").Value = 0.02427782607013916
$ws.Range("C16").Value = 31.48388074484767
$ws.Range("D16").Value = 0.4638807448476712
$ws.Range("E16").Value = 0.2151853454404302
$ws.Range("C17").Value = 31.640887792942
$ws.Range("D17").Value = 0.5208877929419984
$ws.Range("E17").Value = 0.2713240928359862
$ws.Range("C18").Value = 32.12551089741401
$ws.Range("D18").Value = 0.8455108974140089
$ws.Range("E18").Value = 0.7148886776458426
$ws.Range("C19").Value = 31.70430067352617
$ws.Range("D19").Value = 0.3243006735261744
$ws.Range("E19").Value = 0.1051709268495304
$ws.Range("C20").Value = 32.01831504160815
$ws.Range("D20").Value = 0.4383150416081492
$ws.Range("E20").Value = 0.1921200756999535
$ws.Range("C21").Value = 31.6557620332591
$ws.Range("D21").Value = 0.005762033259099297
$ws.Range("E21").Value = 0.00003320102727896647
$ws.Range("C22").Value = 32.55625313419729
$ws.Range("D22").Value = 0.6762531341972924
$ws.Range("E22").Value = 0.4573183015116611
$ws.Range("C23").Value = 32.54997491613913
$ws.Range("D23").Value = 0.2699749161391267
$ws.Range("E23").Value = 0.0728864553443285
$ws.Range("C24").Value = 32.25129507345323
$ws.Range("D24").Value = -0.1987049265467746
$ws.Range("E24").Value = 0.0394836478339591
$ws.Range("C25").Value = 33.63366098406425
$ws.Range("D25").Value = 0.7836609840642481
$ws.Range("E25").Value = 0.6141245379445457
$ws.Range("C26").Value = 32.68773393513963
$ws.Range("D26").Value = -0.2122660648603656
$ws.Range("E26").Value = 0.04505688229130493
$ws.Range("C27").Value = 32.92426366070559
$ws.Range("D27").Value = -0.17573633929441
$ws.Range("E27").Value = 0.0308832609486
$ws.Range("C28").Value = 33.38485988760078
$ws.Range("D28").Value = -0.01514011239921587
$ws.Range("E28").Value = 0.00022922300346089
$ws.Range("C29").Value = 33.76704250188487
$ws.Range("D29").Value = 0.06704250188487038
$ws.Range("E29").Value = 0.004494697058982848
$ws.Range("C30").Value = 34.42621061011622
$ws.Range("D30").Value = 0.3262106101162203
$ws.Range("E30").Value = 0.1064133621523967
$ws.Range("C31").Value = 34.48575269814877
$ws.Range("D31").Value = 0.08575269814877373
$ws.Range("E31").Value = 0.007353525239794702
$ws.Range("C32").Value = 35.01197026628174
$ws.Range("D32").Value = 0.1119702662817446
$ws.Range("E32").Value = 0.01253734053120478
$ws.Range("C33").Value = 35.0747884907683
$ws.Range("D33").Value = -0.2252115092316984
$ws.Range("E33").Value = 0.05072022389041938
$ws.Range("C34").Value = 35.43834716810733
$ws.Range("D34").Value = -0.2616528318926683
$ws.Range("E34").Value = 0.06846220443745293
$ws.Range("C35").Value = 35.65069800646503
$ws.Range("D35").Value = -0.6493019935349693
$ws.Range("E35").Value = 0.4215930788084853
$ws.Range("C36").Value = 36.08379325786624
$ws.Range("D36").Value = -0.7162067421337568
$ws.Range("E36").Value = 0.5129520974778496
$ws.Range("C37").Value = 36.82042427363491
$ws.Range("D37").Value = -0.4795757263650913
$ws.Range("E37").Value = 0.2299928773186049
$ws.Range("C38").Value = 38.2022152009615
$ws.Range("D38").Value = 0.3022152009615056
$ws.Range("E38").Value = 0.0913340276922032
$ws.Range("C39").Value = 38.54488735553488
$ws.Range("D39").Value = 0.04488735553487544
$ws.Range("E39").Value = 0.002014874686914313
$ws.Range("C40").Value = 38.9232925618563
$ws.Range("D40").Value = 0.02329256185630157
$ws.Range("E40").Value = 0.0005425434378296351
$ws.Range("C41").Value = 39.18492527594462
$ws.Range("D41").Value = -0.215074724055377
$ws.Range("E41").Value = 0.04625713692749655
$ws.Range("C42").Value = 40.48918849674822
$ws.Range("D42").Value = 0.5891884967482213
$ws.Range("E42").Value = 0.3471430847004288
$ws.Range("C43").Value = 40.07497914666638
$ws.Range("D43").Value = -0.02502085333362203
$ws.Range("E43").Value = 0.0006260431015426245
$ws.Range("C44").Value = 41.00696585928799
$ws.Range("D44").Value = 0.4069658592879861
$ws.Range("E44").Value = 0.1656212106260089
$ws.Range("C45").Value = 41.00589054627245
$ws.Range("D45").Value = 0.1058905462724482
$ws.Range("E45").Value = 0.01121280778987749
$ws.Range("C46").Value = 41.96942117783217
$ws.Range("D46").Value = 0.7694211778321716
$ws.Range("E46").Value = 0.5920089488966461
$ws.Range("C47").Value = 41.29084426028793
$ws.Range("D47").Value = -0.2091557397120667
$ws.Range("E47").Value = 0.04374612345450179
$ws.Range("C48").Value = 41.11297919396689
$ws.Range("D48").Value = -0.6870208060331038
$ws.Range("E48").Value = 0.4719975879223756
$ws.Range("C49").Value = 42.04974543499112
$ws.Range("D49").Value = -0.1502545650088862
$ws.Range("E49").Value = 0.02257643430600962
$ws.Range("C50").Value = 42.78179685459305
$ws.Range("D50").Value = 0.08179685459304409
$ws.Range("E50").Value = 0.006690725421315598
$ws.Range("C51").Value = 43.32596037961508
$ws.Range("D51").Value = -0.3740396203849272
$ws.Range("E51").Value = 0.1399056376177005

# TOTAL (sum of DELTA), SSE (sum of DELTA^2), and MSE (SSE / count) rows
$ws.Range("C52").Value = 0.3290915639538596
$ws.Range("E52").Value = 8.384676800139479
$ws.Range("E53").Value = 0.1676935360027896
